# Updated cryptos list on Mon Jun  3 09:28:16 UTC 2024 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) figures, and re-rank the few coins
# whose relative order changed (B/C/D/E updated for those rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.944.84'
$ws.Range("E2").Value = '  +1.88%  '

$ws.Range("D3").Value = '3.815.34'
$ws.Range("E3").Value = '  +0.69%  '

$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '626.64'
$ws.Range("E5").Value = '  +4.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.90'
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").Value = '3.814.39'
$ws.Range("E7").Value = '  +0.74%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("E9").Value = '  +0.73%  '

$ws.Range("E10").Value = '  +1.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.454'
$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.93'
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").Value = '4.456.32'

$ws.Range("D16").Value = '3.797.96'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").Value = '68.954.94'
$ws.Range("E17").Value = '  +1.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.02'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.11'
$ws.Range("E19").Value = '  +1.21%  '

$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '464.64'
$ws.Range("E21").Value = '  +0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.65'
$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("E23").Value = '  +1.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000151'
$ws.Range("E24").Value = '  +4.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.68'
$ws.Range("E25").Value = '  +1.27%  '

# Row 26: InternetComputer(DFINITY) -> Fetch.AI
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("E26").Value = '  +2.58%  '

# Row 27: Fetch.AI -> InternetComputer(DFINITY)
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.90'
$ws.Range("E27").Value = '  -0.78%  '

# Row 28: RenderToken -> Dai
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.03%  '

# Row 29: Dai -> RenderToken
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("D30").Value = '3.968.31'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.69'
$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.21'
$ws.Range("E32").Value = '  +0.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.26'
$ws.Range("E33").Value = '  -2.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.11'
$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.06'
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("E37").Value = '  +2.67%  '

$ws.Range("E38").Value = '  +7.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.35'
$ws.Range("E39").Value = '  +4.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.90'
$ws.Range("E40").Value = '  +2.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.979'
$ws.Range("E41").Value = '  -0.64%  '

$ws.Range("E42").Value = '  +0.04%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '157.06'
$ws.Range("E44").Value = '  +3.61%  '

$ws.Range("E45").Value = '  +5.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.298'
$ws.Range("E46").Value = '  +0.48%  '

$ws.Range("E47").Value = '  -1.69%  '

# Row 48: Cosmos -> Arweave
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.54'
$ws.Range("E48").Value = '  -3.18%  '

# Row 49: Stacks -> Cosmos
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.44'
$ws.Range("E49").Value = '  +1.46%  '

# Row 50: Arweave -> Stacks
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.89'
$ws.Range("E50").Value = '  +2.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000279'
$ws.Range("E51").Value = '  +13.75%  '

